$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Row 2: OrchestratorQueueName value (B2) "ProcessABCQueue" -> cleared
$settings.Range("B2").ClearContents()

# Row 7: "ACME_Credential" -> "ACME_LoginCredential" (Name + Value columns)
$settings.Range("A7").Value = "ACME_LoginCredential"
$settings.Range("B7").Value = "ACME_LoginCredential"

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")

# New row 2: ACME_URL asset entry
$assets.Range("A2").Value = "ACME_URL"
$assets.Range("B2").Value = "ACME_URL"
$assets.Range("C2").Value = "REF_ACME"
$assets.Range("D2").Value = "ACME website URL"

# Page setup for the Assets sheet
$ps = $assets.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Row 5: logF_BusinessProcessName value (B5) "Framework" -> "ACME_CalculateClientSecurityHash"
$settings.Range("B5").Value = "ACME_CalculateClientSecurityHash"

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# Row 2: MaxRetryNumber value (B2) 0 -> 1
$constants.Range("B2").Value = 1

# ---------------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------------
$settings.Range("C7").Select()
$assets.Range("D10").Select()
$constants.Activate()
